# Updates england_premier-league_2023-2024 sheet:
#  - Fix 4 row pairs whose match details (columns F:V) had been swapped
#    by mistake: (51,53), (83,84), (86,87), (105,106)
#  - Append newly scraped match row 112 (Wolves 2-1 Tottenham)
#
# NOTE: PowerShell parses "F$r1:V$r1" as the scope-qualified variable
# "$r1:V" (like $env:VAR), silently truncating the string to "F51".
# Always build range addresses with "$($r1)" or string concatenation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $addrA = "F" + $rowA + ":V" + $rowA
    $addrB = "F" + $rowB + ":V" + $rowB

    $valsA = $ws.Range($addrA).Value2
    $valsB = $ws.Range($addrB).Value2

    $ws.Range($addrA).Value2 = $valsB
    $ws.Range($addrB).Value2 = $valsA
}

Swap-Rows 51 53
Swap-Rows 83 84
Swap-Rows 86 87
Swap-Rows 105 106

# Append the new match as row 112, cloning row 111's formatting
# (bold/bordered index style in A, datetime format in E) first.
$ws.Range("A111:V111").Copy($ws.Range("A112:V112"))

$ws.Range("A112").Value2 = 111
$ws.Range("B112").Value2 = "england"
$ws.Range("C112").Value2 = "premier-league"
$ws.Range("D112").Value2 = "2023-2024"
$ws.Range("E112").Value2 = 45241.5625
$ws.Range("F112").Value2 = "Wolves"
$ws.Range("G112").Value2 = 2
$ws.Range("H112").Value2 = "Tottenham"
$ws.Range("I112").Value2 = 1
$ws.Range("J112").Value2 = 3.8
$ws.Range("K112").Value2 = "28/10/2023 20:02"
$ws.Range("L112").Value2 = 3.13
$ws.Range("M112").Value2 = "11/11/2023 13:29"
$ws.Range("N112").Value2 = 3.75
$ws.Range("O112").Value2 = "28/10/2023 20:02"
$ws.Range("P112").Value2 = 3.62
$ws.Range("Q112").Value2 = "11/11/2023 13:29"
$ws.Range("R112").Value2 = 1.95
$ws.Range("S112").Value2 = "28/10/2023 20:02"
$ws.Range("T112").Value2 = 2.31
$ws.Range("U112").Value2 = "11/11/2023 13:29"
$ws.Range("V112").Value2 = "https://www.betexplorer.com/football/england/premier-league/wolves-tottenham/CEzm59UG/"
